# Scheduled-runner style refresh of the per-sheet profit/price columns
# (H:N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW profit sheets. Values are
# plain numeric snapshots (no formulas in this workbook), so each cell is
# just re-stamped with its latest computed value; a few rows gain/lose a
# trailing profit cell (M/N) as the underlying market data changes sign.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 329.5
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H61").Value = 10000
$ws.Range("I61").Value = 10000
$ws.Range("K61").Value = 30000
$ws.Range("M61").Value = -29828
$ws.Range("H74").Value = 14224.0625
$ws.Range("I74").Value = 14505.667
$ws.Range("K74").Value = 14505.667
$ws.Range("M74").Value = -13569.667
$ws.Range("H77").Value = 14224.0625
$ws.Range("I77").Value = 14505.667
$ws.Range("K77").Value = 72528.33499999999
$ws.Range("M77").Value = -67848.33499999999
$ws.Range("H88").Value = 8338157.5
$ws.Range("I88").Value = 25003500
$ws.Range("K88").Value = 25003500
$ws.Range("M88").Value = -25003094
$ws.Range("H91").Value = 8338157.5
$ws.Range("I91").Value = 25003500
$ws.Range("K91").Value = 25003500
$ws.Range("M91").Value = -25002096
$ws.Range("H96").Value = 935.5
$ws.Range("J96").Value = 922.6667
$ws.Range("L96").Value = 2768.0001
$ws.Range("N96").Value = -5514.0001
$ws.Range("H100").Value = 1961.3334
$ws.Range("I100").Value = 1961.3334
$ws.Range("K100").Value = 1961.3334
$ws.Range("M100").Value = -1420.3334
$ws.Range("H113").Value = 5379.5
$ws.Range("I113").Value = 3500
$ws.Range("J113").Value = 6006
$ws.Range("K113").Value = 3500
$ws.Range("L113").Value = 6006
$ws.Range("M113").Value = -246
$ws.Range("N113").Value = -12514
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2293.3076
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H102").Value = 2784.7058
$ws.Range("I102").Value = 2257.077
$ws.Range("K102").Value = 2257.077
$ws.Range("M102").Value = -635.0770000000002
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2970.9167
$ws.Range("I86").Value = 2720.2
$ws.Range("J86").Value = 4224.5
$ws.Range("K86").Value = 2720.2
$ws.Range("L86").Value = 4224.5
$ws.Range("M86").Value = -1597.2
$ws.Range("N86").Value = -6470.5
$ws.Range("H89").Value = 2970.9167
$ws.Range("I89").Value = 2720.2
$ws.Range("J89").Value = 4224.5
$ws.Range("K89").Value = 13601
$ws.Range("L89").Value = 21122.5
$ws.Range("M89").Value = -7985
$ws.Range("N89").Value = -32354.5
$ws.Range("H94").Value = 80004620
$ws.Range("I94").Value = 105268830
$ws.Range("K94").Value = 105268830
$ws.Range("M94").Value = -105268379
$ws.Range("H105").Value = 11306128
$ws.Range("I105").Value = 834770.25
$ws.Range("J105").Value = 22729428
$ws.Range("K105").Value = 834770.25
$ws.Range("L105").Value = 22729428
$ws.Range("M105").Value = -833023.25
$ws.Range("N105").Value = -22732922
$ws.Range("H107").Value = 15387713
$ws.Range("I107").Value = 25644522
$ws.Range("K107").Value = 25644522
$ws.Range("M107").Value = -25642602
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 39999.5
$ws.Range("J28").Value = 39999.5
$ws.Range("L28").Value = 39999.5
$ws.Range("N28").Value = -40489.5
$ws.Range("H99").Value = 2440
$ws.Range("I99").Value = 2000
$ws.Range("K99").Value = 2000
$ws.Range("M99").Value = -502
$ws.Range("H105").Value = 1429.9231
$ws.Range("I105").Value = 1798.4286
$ws.Range("K105").Value = 1798.4286
$ws.Range("M105").Value = -51.42859999999996
$ws.Range("H126").Value = 2440
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 322
$ws.Range("J33").Value = 397
$ws.Range("L33").Value = 2382
$ws.Range("N33").Value = -2948
$ws.Range("H122").Value = 1134.5834
$ws.Range("I122").Value = 493.25
$ws.Range("J122").Value = 1262.85
$ws.Range("K122").Value = 4439.25
$ws.Range("L122").Value = 11365.65
$ws.Range("M122").Value = -1989.25
$ws.Range("N122").Value = -16265.65
$ws.Range("H129").Value = 1837.125
$ws.Range("J129").Value = 2308
$ws.Range("L129").Value = 6924
$ws.Range("N129").Value = -16924
$ws.Range("H131").Value = 11365576
$ws.Range("J131").Value = 1825.2142
$ws.Range("L131").Value = 5475.642599999999
$ws.Range("N131").Value = -15555.6426
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 373447.66
$ws.Range("J95").Value = 373447.66
$ws.Range("L95").Value = 373447.66
$ws.Range("N95").Value = -378939.66
$ws.Range("H98").Value = 24128.2
$ws.Range("J98").Value = 24128.2
$ws.Range("L98").Value = 24128.2
$ws.Range("N98").Value = -30118.2
$ws.Range("H132").Value = 2221.125
$ws.Range("I132").Value = 2041.8422
$ws.Range("K132").Value = 6125.5266
$ws.Range("M132").Value = -3595.5266
$ws.Range("H134").Value = 42800
$ws.Range("J134").Value = 42800
$ws.Range("L134").Value = 128400
$ws.Range("N134").Value = -133470
$ws.Range("H136").Value = 56729.4
$ws.Range("J136").Value = 56729.4
$ws.Range("L136").Value = 170188.2
$ws.Range("N136").Value = -175288.2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1762.4286
$ws.Range("I16").Value = 1700
$ws.Range("J16").Value = 1845.6666
$ws.Range("K16").Value = 1700
$ws.Range("L16").Value = 1845.6666
$ws.Range("M16").Value = -1530
$ws.Range("N16").Value = -2185.6666
$ws.Range("H46").Value = 3783.1667
$ws.Range("I46").Value = 1925
$ws.Range("J46").Value = 7499.5
$ws.Range("K46").Value = 1925
$ws.Range("L46").Value = 7499.5
$ws.Range("M46").Value = -1737
$ws.Range("N46").Value = -7875.5
$ws.Range("H61").Value = 900
$ws.Range("I61").Value = 900
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 900
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -698
$ws.Range("N61").ClearContents()
$ws.Range("H100").Value = 5961.231
$ws.Range("I100").Value = 5645.4546
$ws.Range("J100").Value = 7698
$ws.Range("K100").Value = 5645.4546
$ws.Range("L100").Value = 7698
$ws.Range("M100").Value = -5104.4546
$ws.Range("N100").Value = -8780
$ws.Range("H113").Value = 900
$ws.Range("I113").Value = 900
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 900
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1270
$ws.Range("N113").ClearContents()
$ws.Range("H136").Value = 3379.9412
$ws.Range("I136").Value = 2389.9285
$ws.Range("K136").Value = 7169.7855
$ws.Range("M136").Value = -4619.7855
